# TestReport_Group04.xlsx - "change TestPlan + Test Report summary"
#
# Adds the Test Plan block (function list + counts) to the top of the
# "Test summary report" sheet, and flips which sheet/cell is active so the
# workbook re-opens on "Test summary report"!A3 instead of "Bug report"!C3.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test summary report")
$ws2 = $wb.Worksheets.Item("Bug report")

# --- New Test Plan rows under the existing Tester/Date rows -----------------
$ws1.Range("A3").Value  = "Function be Tested"
$ws1.Range("A4").Value  = " number of test cases"
$ws1.Range("A5").Value  = "number of passed test cases,"
$ws1.Range("A6").Value  = "number of failed test cases of a function"
$ws1.Range("A7").Value  = "Function 01: Load Data and Create a List Of Object in Scroll View"
$ws1.Range("A8").Value  = "Function 02: Switch Between Lists"
$ws1.Range("A9").Value  = "Function 03: Click Item and Load Buy View"
$ws1.Range("A10").Value = "Function 04: Load Item Information in Buy View"
$ws1.Range("A11").Value = "Function 05: Quantity View"
$ws1.Range("A12").Value = "Function 06: Facebook Sharing"
$ws1.Range("A13").Value = "Function 07: Progress Tracking"

# --- Column widths to fit the new text (authored widths, minus the COM
#     ColumnWidth -> stored-width padding baked into this host) ------------
$ws1.Columns.Item(1).ColumnWidth = 57.592447916666664
$ws1.Columns.Item(2).ColumnWidth = 13.451822916666666
$ws1.Columns.Item(3).ColumnWidth = 27.877604166666668

# --- Selection / active-sheet swap: "Bug report" loses focus, "Test
#     summary report" gains it with A3 selected -----------------------------
$ws2.Range("B5").Select()
$ws1.Activate()
$ws1.Range("A3").Select()
